$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header (AUC) and column width for the new column H
$ws.Range("H1").Value = "AUC"

# Narrow column C (was 14.7109375 chars, now 6.7109375 chars)
$ws.Columns.Item(3).ColumnWidth = 5.8

# Give the new column H the same width as columns D-G (12.7109375 chars)
$ws.Columns.Item(8).ColumnWidth = 11.8

# Updated metric values for rows 2-31 (new bayes results)
$ws.Range("C2").Value = 2699.5999999999999
$ws.Range("D2").Value = 0.64737740533829924
$ws.Range("E2").Value = 0.93549313358302122
$ws.Range("F2").Value = 0.35923845193508119
$ws.Range("G2").Value = 0.50352243136166175
$ws.Range("H2").Value = 0.64736579275905126
$ws.Range("C3").Value = 80.400000000000006
$ws.Range("D3").Value = 0.91210117939168212
$ws.Range("E3").Value = 0.86762796504369555
$ws.Range("F3").Value = 0.956629213483146
$ws.Range("G3").Value = 0.91599312731606375
$ws.Range("H3").Value = 0.91212858926342089
$ws.Range("C4").Value = 2039.4000000000001
$ws.Range("D4").Value = 0.68629112352576049
$ws.Range("E4").Value = 0.93549313358302122
$ws.Range("F4").Value = 0.43704119850187262
$ws.Range("G4").Value = 0.57996404383954947
$ws.Range("H4").Value = 0.68626716604244709
$ws.Range("C5").Value = 2591.8000000000002
$ws.Range("D5").Value = 0.65517690875232781
$ws.Range("E5").Value = 0.93769038701622986
$ws.Range("F5").Value = 0.37267166042446942
$ws.Range("G5").Value = 0.51869582032522388
$ws.Range("H5").Value = 0.65518102372034959
$ws.Range("C6").Value = 2745.8000000000002
$ws.Range("D6").Value = 0.65180633147113609
$ws.Range("E6").Value = 0.94220973782771544
$ws.Range("F6").Value = 0.36144818976279652
$ws.Range("G6").Value = 0.50847230149375533
$ws.Range("H6").Value = 0.6518289637952559
$ws.Range("C7").Value = 1740.0999999999999
$ws.Range("D7").Value = 0.64277777777777789
$ws.Range("E7").Value = 0.87333333333333329
$ws.Range("F7").Value = 0.41222222222222216
$ws.Range("G7").Value = 0.53567103225659896
$ws.Range("H7").Value = 0.64277777777777789
$ws.Range("C8").Value = 659.79999999999995
$ws.Range("D8").Value = 0.78421787709497204
$ws.Range("E8").Value = 0.92439450686641711
$ws.Range("F8").Value = 0.64408239700374526
$ws.Range("G8").Value = 0.74814851458377751
$ws.Range("H8").Value = 0.78423845193508113
$ws.Range("C9").Value = 994.60000000000002
$ws.Range("D9").Value = 0.76973929236499061
$ws.Range("E9").Value = 0.94323345817727833
$ws.Range("F9").Value = 0.59625468164794015
$ws.Range("G9").Value = 0.72037182119388943
$ws.Range("H9").Value = 0.76974406991260913
$ws.Range("C10").Value = 1394.3
$ws.Range("D10").Value = 0.71134388578522667
$ws.Range("E10").Value = 0.91774032459425725
$ws.Range("F10").Value = 0.50495630461922603
$ws.Range("G10").Value = 0.63529941905750453
$ws.Range("H10").Value = 0.71134831460674164
$ws.Range("C11").Value = 2477.4000000000001
$ws.Range("D11").Value = 0.60843575418994411
$ws.Range("E11").Value = 0.88428214731585508
$ws.Range("F11").Value = 0.33264669163545568
$ws.Range("G11").Value = 0.45903742924697344
$ws.Range("H11").Value = 0.60846441947565544
$ws.Range("C12").Value = 4519.3999999999996
$ws.Range("D12").Value = 0.55506207324643075
$ws.Range("E12").Value = 0.92771535580524345
$ws.Range("F12").Value = 0.18247191011235955
$ws.Range("G12").Value = 0.28833038339741912
$ws.Range("H12").Value = 0.55509363295880154
$ws.Range("C13").Value = 94.099999999999994
$ws.Range("D13").Value = 0.9115642458100558
$ws.Range("E13").Value = 0.86205992509363294
$ws.Range("F13").Value = 0.96108614232209744
$ws.Range("G13").Value = 0.9157619204511549
$ws.Range("H13").Value = 0.91157303370786524
$ws.Range("C14").Value = 2738.9000000000001
$ws.Range("D14").Value = 0.61678150217256356
$ws.Range("E14").Value = 0.906541822721598
$ws.Range("F14").Value = 0.32697877652933832
$ws.Range("G14").Value = 0.45960931754896917
$ws.Range("H14").Value = 0.61676029962546808
$ws.Range("C15").Value = 2379.6999999999998
$ws.Range("D15").Value = 0.66131905648665434
$ws.Range("E15").Value = 0.92993757802746568
$ws.Range("F15").Value = 0.39275905118601739
$ws.Range("G15").Value = 0.53382979881543136
$ws.Range("H15").Value = 0.66134831460674159
$ws.Range("C16").Value = 2258.3000000000002
$ws.Range("D16").Value = 0.67019553072625704
$ws.Range("E16").Value = 0.9332958801498128
$ws.Range("F16").Value = 0.40712858926342071
$ws.Range("G16").Value = 0.55129643934317085
$ws.Range("H16").Value = 0.67021223470661684
$ws.Range("C17").Value = 3094.5
$ws.Range("D17").Value = 0.59833333333333327
$ws.Range("E17").Value = 0.90444444444444438
$ws.Range("F17").Value = 0.29222222222222227
$ws.Range("G17").Value = 0.41751705834148878
$ws.Range("H17").Value = 0.59833333333333327
$ws.Range("C18").Value = 4338.1999999999998
$ws.Range("D18").Value = 0.60844506517690866
$ws.Range("E18").Value = 0.97439450686641715
$ws.Range("F18").Value = 0.24247191011235958
$ws.Range("G18").Value = 0.38141999019705369
$ws.Range("H18").Value = 0.6084332084893882
$ws.Range("C19").Value = 549.29999999999995
$ws.Range("D19").Value = 0.83815021725636252
$ws.Range("E19").Value = 0.96554307116104887
$ws.Range("F19").Value = 0.71071161048689147
$ws.Range("G19").Value = 0.81384987256700703
$ws.Range("H19").Value = 0.83812734082397
$ws.Range("C20").Value = 1614
$ws.Range("D20").Value = 0.71967411545623849
$ws.Range("E20").Value = 0.93993757802746569
$ws.Range("F20").Value = 0.49947565543071165
$ws.Range("G20").Value = 0.6386168840616967
$ws.Range("H20").Value = 0.71970661672908864
$ws.Range("C21").Value = 3127.3000000000002
$ws.Range("D21").Value = 0.62682184978274358
$ws.Range("E21").Value = 0.93661672908863935
$ws.Range("F21").Value = 0.31695380774032461
$ws.Range("G21").Value = 0.45826849405554781
$ws.Range("H21").Value = 0.6267852684144819
$ws.Range("C22").Value = 3807.6999999999998
$ws.Range("D22").Value = 0.57232153941651143
$ws.Range("E22").Value = 0.91440699126092395
$ws.Range("F22").Value = 0.23028714107365794
$ws.Range("G22").Value = 0.34870695852499489
$ws.Range("H22").Value = 0.57234706616729081
$ws.Range("C23").Value = 131.90000000000001
$ws.Range("D23").Value = 0.87817504655493495
$ws.Range("E23").Value = 0.81978776529338315
$ws.Range("F23").Value = 0.93654182272159792
$ws.Range("G23").Value = 0.88473300809197486
$ws.Range("H23").Value = 0.87816479400749059
$ws.Range("C24").Value = 2094.9000000000001
$ws.Range("D24").Value = 0.65350403476101793
$ws.Range("E24").Value = 0.90660424469413248
$ws.Range("F24").Value = 0.40037453183520599
$ws.Range("G24").Value = 0.53457705261663979
$ws.Range("H24").Value = 0.65348938826466918
$ws.Range("C25").Value = 2412.8000000000002
$ws.Range("D25").Value = 0.66959342023587842
$ws.Range("E25").Value = 0.94212234706616726
$ws.Range("F25").Value = 0.39711610486891386
$ws.Range("G25").Value = 0.5448173463924032
$ws.Range("H25").Value = 0.66961922596754064
$ws.Range("C26").Value = 2084.5999999999999
$ws.Range("D26").Value = 0.66631284916201117
$ws.Range("E26").Value = 0.91878901373283406
$ws.Range("F26").Value = 0.41379525593008742
$ws.Range("G26").Value = 0.5519961896031772
$ws.Range("H26").Value = 0.66629213483146077
$ws.Range("C27").Value = 2946
$ws.Range("D27").Value = 0.57666666666666655
$ws.Range("E27").Value = 0.87777777777777788
$ws.Range("F27").Value = 0.27555555555555555
$ws.Range("G27").Value = 0.39345901727968574
$ws.Range("H27").Value = 0.57666666666666655
$ws.Range("C28").Value = 1710
$ws.Range("D28").Value = 0.6907697082557418
$ws.Range("E28").Value = 0.91097378277153562
$ws.Range("F28").Value = 0.47049937578027468
$ws.Range("G28").Value = 0.59601714987287135
$ws.Range("H28").Value = 0.69073657927590515
$ws.Range("C29").Value = 880.5
$ws.Range("D29").Value = 0.79813469894475486
$ws.Range("E29").Value = 0.96108614232209733
$ws.Range("F29").Value = 0.63516853932584272
$ws.Range("G29").Value = 0.75810812728863097
$ws.Range("H29").Value = 0.79812734082397008
$ws.Range("C30").Value = 1746.8
$ws.Range("D30").Value = 0.70970515207945373
$ws.Range("E30").Value = 0.94109862671660438
$ws.Range("F30").Value = 0.47836454431960052
$ws.Range("G30").Value = 0.62235015793752813
$ws.Range("H30").Value = 0.70973158551810245
$ws.Range("C31").Value = 2760.5999999999999
$ws.Range("D31").Value = 0.62179702048417129
$ws.Range("E31").Value = 0.91323345817727852
$ws.Range("F31").Value = 0.33033707865168538
$ws.Range("G31").Value = 0.46489525789444242
$ws.Range("H31").Value = 0.62178526841448201
